# feat: add 2022-Q1 data
# Inserts a new "2022-Q1" worksheet (between "2021-Q4" and "总计") with
# per-fund holding detail, and updates the "总计" summary sheet with a new
# leading row for the 2022-Q1 quarter (shifting the 2021-Q4 row down).

$wb = $excel.ActiveWorkbook

$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsTotal = $wb.Worksheets.Item("总计")

# Style source cells that already carry the workbook's "header / index
# column" look (bold, boxed, centered) - reuse their exact style via
# copy/paste-special instead of re-deriving font+border+alignment by hand.
$headerStyleSrc = $wsQ4.Range("B1")
$indexStyleSrc = $wsQ4.Range("A2")

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted right before "总计"
# ---------------------------------------------------------------------
$wsQ1 = $wb.Worksheets.Add($wsTotal)
$wsQ1.Name = "2022-Q1"

# Inserting a sheet shifts positional handles, so re-resolve "总计" by name
# rather than relying on the (now stale) $wsTotal reference captured above.
$wsTotal = $wb.Worksheets.Item("总计")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $wsQ1.Cells.Item(1, 2 + $i)
    $headerStyleSrc.Copy()
    $cell.PasteSpecial(-4122)
    $cell.Value = $headers[$i]
}

# index, code, name, fund size, stock position, position ratio, market value, rank
$rows = @(
    @(0,  "630010", "华商价值精选混合", "4.93", "89.49", "3.14", "0.1548", 8),
    @(1,  "002938", "中银证券健康产业灵活配置混合", "1.82", "59.84", "4.59", "0.0835", 2),
    @(2,  "011269", "中银证券优势制造股票型证券投资基金A", "1.39", "93.51", "5.74", "0.0798", 4),
    @(3,  "005108", "圆信永丰双利优选定期开放灵活配置混合", "1.89", "94.60", "3.83", "0.0724", 10),
    @(4,  "970032", "东海证券海睿进取灵活配置混合型集合资产管理计划A", "1.57", "78.62", "3.35", "0.0526", 6),
    @(5,  "160921", "大成多策略混合(LOF)", "1.13", "79.19", "3.89", "0.0440", 9),
    @(6,  "970033", "东海证券海睿进取灵活配置混合型集合资产管理计划B", "1.18", "78.62", "3.35", "0.0395", 6),
    @(7,  "001965", "圆信永丰兴源灵活配置混合A", "0.76", "93.43", "4.35", "0.0331", 10),
    @(8,  "010434", "红土创新医疗保健股票", "0.75", "92.96", "4.16", "0.0312", 10),
    @(9,  "630006", "华商产业升级混合", "0.98", "87.95", "3.14", "0.0308", 8),
    @(10, "001421", "南方量化成长股票", "1.70", "92.11", "1.78", "0.0303", 10),
    @(11, "011270", "中银证券优势制造股票型证券投资基金C", "0.21", "93.51", "5.74", "0.0121", 4),
    @(12, "001966", "圆信永丰兴源灵活配置混合C", "0.25", "93.43", "4.35", "0.0109", 10),
    @(13, "006274", "圆信永丰医药健康混合", "0.18", "93.66", "4.15", "0.0075", 9),
    @(14, "005536", "渤海汇金量化成长混合", "0.61", "88.57", "1.20", "0.0073", 2),
    @(15, "970050", "东海海睿锐意3个月定开灵活配置混合", "0.17", "78.42", "3.19", "0.0054", 6),
    @(16, "002952", "建信多因子量化股票", "0.10", "91.47", "2.79", "0.0028", 7),
    @(17, "007808", "北信瑞丰量化优选灵活配置混合", "0.24", "79.84", "1.04", "0.0025", 10)
)

foreach ($r in $rows) {
    $rowNum = $r[0] + 2

    $aCell = $wsQ1.Cells.Item($rowNum, 1)
    $indexStyleSrc.Copy()
    $aCell.PasteSpecial(-4122)
    $aCell.Value = $r[0]

    # text-valued columns (kept as text even though they look numeric)
    foreach ($pair in @(@(2, $r[1]), @(3, $r[2]), @(4, $r[3]), @(5, $r[4]), @(6, $r[5]), @(7, $r[6]))) {
        $col = $pair[0]
        $val = $pair[1]
        $cell = $wsQ1.Cells.Item($rowNum, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = "Normal"
    }

    # rank column stays numeric
    $wsQ1.Cells.Item($rowNum, 8).Value = $r[7]
}

# ---------------------------------------------------------------------
# 2) Update "总计": insert a new row 2 for 2022-Q1, pushing 2021-Q4 to row 3
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()

$a2 = $wsTotal.Cells.Item(2, 1)
$indexStyleSrc.Copy()
$a2.PasteSpecial(-4122)
$a2.Value = 0

$b2 = $wsTotal.Cells.Item(2, 2)
$b2.NumberFormat = "@"
$b2.Value = "2022-Q1"
$b2.Style = "Normal"

$wsTotal.Cells.Item(2, 3).Style = "Normal"
$wsTotal.Cells.Item(2, 3).Value = 18
$wsTotal.Cells.Item(2, 4).Style = "Normal"
$wsTotal.Cells.Item(2, 4).Value = 0.7

# The row that was pushed down (now row 3, originally "2021-Q4" at index 0)
# keeps its running index in sync with its new position, and needs its
# style restored (Insert() drags a copy of the row-2 formatting down).
$a3 = $wsTotal.Cells.Item(3, 1)
$indexStyleSrc.Copy()
$a3.PasteSpecial(-4122)
$a3.Value = 1

$wsTotal.Cells.Item(3, 3).Style = "Normal"
$wsTotal.Cells.Item(3, 4).Style = "Normal"
